$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '27.379.13'
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = '  -1.58%  '

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.655.15'
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = '  -0.41%  '

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = '  -0.15%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '213.22'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  -0.55%  '

# Row 6
$ws.Range("E6").Value = '  -0.27%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.00'
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = '  -0.07%  '

# Row 8
$ws.Range("E8").Value = '  +1.24%  '

# Row 9
$ws.Range("E9").Value = '  -0.69%  '

# Row 10
$ws.Range("E10").Value = '  -1.22%  '

# Row 11
$ws.Range("E11").Value = '  -0.44%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.890.41'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  -0.41%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.661.10'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  -0.08%  '

# Row 14
$ws.Range("E14").Value = '  -1.43%  '

# Row 15
$ws.Range("E15").Value = '  +3.50%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '65.54'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  -0.60%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '27.387.51'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  -1.46%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '231.62'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  -6.45%  '

# Row 19
$ws.Range("E19").Value = '  -0.73%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.46'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  -0.22%  '

# Row 21
$ws.Range("E21").Value = '  -0.02%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.36'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  -2.41%  '

# Row 23
$ws.Range("E23").Value = '  +0.43%  '

# Row 24
$ws.Range("E24").Value = '  -0.88%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '147.72'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  +0.66%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '7.10'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  -1.03%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '15.88'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  -2.24%  '

# Row 28
$ws.Range("E28").Value = '  -0.24%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.112'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  -0.09%  '

# Row 30
$ws.Range("E30").Value = '  -0.39%  '

# Row 31
$ws.Range("E31").Value = '  -4.16%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.29'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  -1.39%  '

# Row 33
$ws.Range("B33").Value = 'InternetComputer(DFINITY)'
$ws.Range("C33").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.13'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  +0.15%  '

# Row 34
$ws.Range("B34").Value = 'Maker'
$ws.Range("C34").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.426.85'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  -0.59%  '

# Row 35
$ws.Range("E35").Value = '  +0.79%  '

# Row 36
$ws.Range("E36").Value = '  -0.78%  '

# Row 37
$ws.Range("E37").Value = '  -2.05%  '

# Row 38
$ws.Range("E38").Value = '  -1.31%  '

# Row 39
$ws.Range("E39").Value = '  +0.22%  '

# Row 40
$ws.Range("E40").Value = '  -0.30%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.00'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  -0.16%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '5.50'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  +1.68%  '

# Row 43
$ws.Range("B43").Value = 'TrustWalletToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.797'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  +1.16%  '

# Row 44
$ws.Range("B44").Value = 'Aave'
$ws.Range("C44").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '64.92'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  -6.30%  '

# Row 45
$ws.Range("B45").Value = 'MXToken'
$ws.Range("C45").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.22'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  -0.08%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.798.24'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  -0.30%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.68'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  -1.09%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '87.96'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  -1.26%  '

# Row 49
$ws.Range("E49").Value = '  -3.01%  '

# Row 50
$ws.Range("E50").Value = '  -0.29%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '7.73'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  -0.99%  '
